# Append two new daily rows (date=45988/2025-11-27 and date=45989/2025-11-28)
# to each of the 4 worksheets, matching the new "remn_amt" values for each
# company. Row 120 (2025-11-28) has a remn_amt of 0 on every sheet.

$wb = $excel.ActiveWorkbook

$sheetValues = @{
    1 = @(490, 0)
    2 = @(3602, 0)
    3 = @(2775, 0)
    4 = @(1261, 0)
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $values = $sheetValues[$i]

    $ws.Cells.Item(119, 1).Value = 45988
    $ws.Cells.Item(119, 2).Value = $values[0]

    $ws.Cells.Item(120, 1).Value = 45989
    $ws.Cells.Item(120, 2).Value = $values[1]

    # Match the date-styled format used by the rest of column A (e.g. A118).
    $ws.Range("A119:A120").NumberFormat = $ws.Range("A118").NumberFormat
}
